$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 12028.1
$ws.Range("I64").Value = 9594.4
$ws.Range("J64").Value = 14461.8
$ws.Range("K64").Value = 9594.4
$ws.Range("L64").Value = 14461.8
$ws.Range("M64").Value = -9346.4
$ws.Range("N64").Value = -14957.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 12028.1
$ws.Range("I67").Value = 9594.4
$ws.Range("J67").Value = 14461.8
$ws.Range("K67").Value = 9594.4
$ws.Range("L67").Value = 14461.8
$ws.Range("M67").Value = -8736.4
$ws.Range("N67").Value = -16177.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 16446.047
$ws.Range("I69").Value = 10996.667
$ws.Range("J69").Value = 17354.277
$ws.Range("K69").Value = 32990.001
$ws.Range("L69").Value = 52062.83099999999
$ws.Range("M69").Value = -32116.001
$ws.Range("N69").Value = -53810.83099999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 16446.047
$ws.Range("I72").Value = 10996.667
$ws.Range("J72").Value = 17354.277
$ws.Range("K72").Value = 98970.003
$ws.Range("L72").Value = 156188.493
$ws.Range("M72").Value = -94602.003
$ws.Range("N72").Value = -164924.493

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 79617.30499999999
$ws.Range("I111").Value = 1003.8889
$ws.Range("K111").Value = 3011.6667
$ws.Range("M111").Value = 55.33329999999978

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 64550.75
$ws.Range("I125").Value = 1004.7778
$ws.Range("K125").Value = 9043.0002
$ws.Range("M125").Value = -6583.0002

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 2045.9375
$ws.Range("J131").Value = 3975
$ws.Range("L131").Value = 11925
$ws.Range("N131").Value = -22005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1202.8334
$ws.Range("J97").Value = 822
$ws.Range("L97").Value = 822
$ws.Range("N97").Value = -1814

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4495.2085
$ws.Range("I132").Value = 2644.3
$ws.Range("J132").Value = 13749.75
$ws.Range("K132").Value = 7932.900000000001
$ws.Range("L132").Value = 41249.25
$ws.Range("M132").Value = -5402.900000000001
$ws.Range("N132").Value = -46309.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 20252
$ws.Range("J80").Value = 14621.857
$ws.Range("L80").Value = 14621.857
$ws.Range("N80").Value = -16617.857

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 20252
$ws.Range("J83").Value = 14621.857
$ws.Range("L83").Value = 73109.285
$ws.Range("N83").Value = -83093.285

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 2863346.8
$ws.Range("I86").Value = 4449296
$ws.Range("K86").Value = 4449296
$ws.Range("M86").Value = -4448173

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 2863346.8
$ws.Range("I89").Value = 4449296
$ws.Range("K89").Value = 22246480
$ws.Range("M89").Value = -22240864

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4919.154
$ws.Range("J99").Value = 3768.5652
$ws.Range("L99").Value = 3768.5652
$ws.Range("N99").Value = -6764.5652

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 4919.154
$ws.Range("J126").Value = 3768.5652
$ws.Range("L126").Value = 11305.6956
$ws.Range("N126").Value = -16245.6956

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 1374.8334
$ws.Range("J97").Value = 1349.875
$ws.Range("L97").Value = 4049.625
$ws.Range("N97").Value = -5041.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 5306613
$ws.Range("I121").Value = 1799.75
$ws.Range("J121").Value = 7959020
$ws.Range("K121").Value = 5399.25
$ws.Range("L121").Value = 23877060
$ws.Range("M121").Value = -4089.25
$ws.Range("N121").Value = -23879680

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 4575.88
$ws.Range("J131").Value = 5707
$ws.Range("L131").Value = 17121
$ws.Range("N131").Value = -27201

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 4004.5557
$ws.Range("I137").Value = 2381.3333
$ws.Range("J137").Value = 7251
$ws.Range("K137").Value = 7143.999899999999
$ws.Range("L137").Value = 21753
$ws.Range("M137").Value = -2043.999899999999
$ws.Range("N137").Value = -31953

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3879.2354
$ws.Range("I102").Value = 2437.9412
$ws.Range("K102").Value = 2437.9412
$ws.Range("M102").Value = -815.9412000000002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 167689.83
$ws.Range("I7").Value = 167689.83
$ws.Range("K7").Value = 167689.83
$ws.Range("M7").Value = -167577.83

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 8213.75
$ws.Range("I68").Value = 5257.8887
$ws.Range("K68").Value = 5257.8887
$ws.Range("M68").Value = -4508.8887

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 8213.75
$ws.Range("I71").Value = 5257.8887
$ws.Range("K71").Value = 26289.4435
$ws.Range("M71").Value = -22545.4435

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("M92").ClearContents()
$ws.Range("N92").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H99").Value = 43094.668
$ws.Range("I99").Value = 29999.5
$ws.Range("K99").Value = 29999.5
$ws.Range("M99").Value = -27004.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 167689.83
$ws.Range("I126").Value = 167689.83
$ws.Range("K126").Value = 503069.49
$ws.Range("M126").Value = -500599.49

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 6479.3013
$ws.Range("I136").Value = 4312.2915
$ws.Range("J136").Value = 10639.96
$ws.Range("K136").Value = 12936.8745
$ws.Range("L136").Value = 31919.88
$ws.Range("M136").Value = -10386.8745
$ws.Range("N136").Value = -37019.88

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4777
$ws.Range("J81").Value = 5666
$ws.Range("L81").Value = 11332
$ws.Range("N81").Value = -13454

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 4777
$ws.Range("J84").Value = 5666
$ws.Range("L84").Value = 56660
$ws.Range("N84").Value = -67268

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 9030.357
$ws.Range("J96").Value = 19999.4
$ws.Range("L96").Value = 19999.4
$ws.Range("N96").Value = -22745.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1672.6
$ws.Range("I126").Value = 1049.9166
$ws.Range("J126").Value = 4163.3335
$ws.Range("K126").Value = 3149.7498
$ws.Range("L126").Value = 12490.0005
$ws.Range("M126").Value = -679.7498000000001
$ws.Range("N126").Value = -17430.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3957.17
$ws.Range("I132").Value = 2995.468
$ws.Range("K132").Value = 8986.403999999999
$ws.Range("M132").Value = -6456.403999999999
